$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26 (shifts existing rows 26-58 down to 27-59,
# carrying their formatting/styles with them).
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly price entry.
$ws.Range("A26").Value = 4
$ws.Range("B26").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C26").Value = "Los Lagos"
$ws.Range("D26").Value = 44973
$ws.Range("E26").Value = 10
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100101
$ws.Range("H26").Value = "Berries"
$ws.Range("I26").Value = 100101001
$ws.Range("J26").Value = "Arándano (blue)"
$ws.Range("K26").Value = "Sin especificar"
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 200
$ws.Range("N26").Value = 2000
$ws.Range("O26").Value = 2200
$ws.Range("P26").Value = 2100
$ws.Range("Q26").Value = "$/bandeja 2 kilos"
$ws.Range("R26").Value = "Provincia de Curicó"
$ws.Range("S26").Value = 1050
$ws.Range("T26").Value = 2
